$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.486292719841003
$ws.Range("B1").Value = 3.591972351074219
$ws.Range("C1").Value = 2.601009368896484
$ws.Range("D1").Value = 1.333489775657654
$ws.Range("E1").Value = 0.7662955522537231
